# Switch Pit and Tank tab names
# ------------------------------------------------------------------
# The workbook has three sheets: Pit, Tank, Borehole. Each of Pit/Tank
# hosts its own XY-scatter chart that plots its own A2:A9 / B2:B9 data.
# This change swaps the tab names "Pit" and "Tank" (the underlying
# sheetId/rId/data for each physical sheet stay put - only the
# displayed tab names trade places), updates the charts' series-name
# formulas so they keep following their own (renamed) host sheet, and
# removes the stray hidden "_xlchart.v1.*" defined names that used to
# back the Borehole chart.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Swap the "Pit" and "Tank" worksheet tab names ---------------
# Use a temporary name so the two renames don't collide.
$pitSheet  = $wb.Worksheets.Item("Pit")
$tankSheet = $wb.Worksheets.Item("Tank")

$pitSheet.Name  = "PitTmpSwapName"
$tankSheet.Name = "Pit"
$pitSheet.Name  = "Tank"

# --- 2. Remove the hidden _xlchart.v1.* defined names ----------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- 3. Re-point each chart's series-name formula at its own sheet ---
# (the sheet that is now named "Tank" is the one that used to be "Pit",
# and vice-versa; each chart keeps showing the data physically on its
# own sheet, just under the new tab name)
$sheetNowTank = $wb.Worksheets.Item("Tank")
$sheetNowPit  = $wb.Worksheets.Item("Pit")

$chartOnTank = $sheetNowTank.ChartObjects().Item(1).Chart
$seriesOnTank = $chartOnTank.SeriesCollection().Item(1)
$seriesOnTank.Name = "=Tank!`$B`$1"

$chartOnPit = $sheetNowPit.ChartObjects().Item(1).Chart
$seriesOnPit = $chartOnPit.SeriesCollection().Item(1)
$seriesOnPit.Name = "=Pit!`$B`$1"

# --- 4. Restore the active-cell selection recorded on the sheet that -
# is now named "Tank" (physically the former "Pit" sheet)
[void]$sheetNowTank.Activate()
[void]$sheetNowTank.Range("G25").Select()
